$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last two data rows (original rows 6 and 7, Sending cluster = MuSCs,
# Target cluster = FAPs / MuSCs) - the new TPM run only has 4 data rows.
$ws.Rows("6:7").Delete()

# Row 2: ECs -> ECs (recalculated TPM figures)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.122548
$ws.Range("H2").Value = 0.367644
$ws.Range("I2").Value = 0.4789089061706202
$ws.Range("J2").Value = 0.4789089061706202
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.051307
$ws.Range("N2").Value = 0.153921
$ws.Range("O2").Value = 0.04260998881052132
$ws.Range("P2").Value = 0.04260998881052132
$ws.Range("Q2").Value = 0.006287570236
$ws.Range("R2").Value = 0.056588132124
$ws.Range("S2").Value = 0.02040630313318913
$ws.Range("T2").Value = 0.02040630313318913

# Row 3: ECs -> FAPs (recalculated TPM figures)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.122548
$ws.Range("H3").Value = 0.367644
$ws.Range("I3").Value = 0.4789089061706202
$ws.Range("J3").Value = 0.4789089061706202
$ws.Range("O3").Value = 0.9573900111894786
$ws.Range("P3").Value = 0.9573900111894786
$ws.Range("Q3").Value = 0.1412733752493333
$ws.Range("R3").Value = 1.271460377244
$ws.Range("S3").Value = 0.4585026030374311
$ws.Range("T3").Value = 0.458502603037431

# Row 4: was ECs -> MuSCs, now MuSCs -> ECs (recalculated TPM figures)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "ECs"
$ws.Range("G4").Value = 0.133342
$ws.Range("H4").Value = 0.400026
$ws.Range("I4").Value = 0.5210910938293798
$ws.Range("J4").Value = 0.5210910938293798
$ws.Range("M4").Value = 0.051307
$ws.Range("N4").Value = 0.153921
$ws.Range("O4").Value = 0.04260998881052132
$ws.Range("P4").Value = 0.04260998881052132
$ws.Range("Q4").Value = 0.006841377994
$ws.Range("R4").Value = 0.061572401946
$ws.Range("S4").Value = 0.02220368567733219
$ws.Range("T4").Value = 0.02220368567733219

# Row 5: was MuSCs -> ECs, now MuSCs -> FAPs (recalculated TPM figures)
$ws.Range("D5").Value = "FAPs"
$ws.Range("G5").Value = 0.133342
$ws.Range("H5").Value = 0.400026
$ws.Range("I5").Value = 0.5210910938293798
$ws.Range("J5").Value = 0.5210910938293798
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.152800333333333
$ws.Range("N5").Value = 3.458401
$ws.Range("O5").Value = 0.9573900111894786
$ws.Range("P5").Value = 0.9573900111894786
$ws.Range("Q5").Value = 0.1537167020473333
$ws.Range("R5").Value = 1.383450318426
$ws.Range("S5").Value = 0.4988874081520476
$ws.Range("T5").Value = 0.4988874081520476
